# Remove the readme.txt hyperlink reference used for "Add the MCR directory
# to the system path (readme.txt)." -- the link was incorrect, so the
# parenthetical "(readme.txt)" (and its hyperlink) is dropped entirely,
# leaving "...to the system path. This step may be unnecessary...".

$d = $word.ActiveDocument

$find = $d.Content.Find
$find.ClearFormatting()
$find.Replacement.ClearFormatting()

$replaced = $find.Execute(
    " (readme.txt).This",  # FindText
    $false,                 # MatchCase
    $false,                 # MatchWholeWord
    $false,                 # MatchWildcards
    $false,                 # MatchSoundsLike
    $false,                 # MatchAllWordForms
    $true,                  # Forward
    1,                      # Wrap (wdFindContinue)
    $false,                 # Format
    ". This",               # ReplaceWith
    2                       # Replace (wdReplaceAll)
)

Write-Output ("readme.txt hyperlink reference removed: " + $replaced)
